$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Merge the "7 devices ... date r" / "ange" split back into a single
# run with the corrected/full text "... every single date range". This also
# removes the original mid-sentence "_GoBack" bookmark that used to sit at
# that split point, since the Find/Replace rewrites the whole paragraph text.
# ---------------------------------------------------------------------------
$mergeRange = $d.Content
$mergeRange.Find.Execute( `
    "7 devices, 6 weather types, 1 month of date range, data of every device for every single date range", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "7 devices, 6 weather types, 1 month of date range, data of every device for every single date range", 2)

# Locate the (now single-run) paragraph that holds this sentence so we can
# anchor the new bookmarks relative to it.
$sevenDevicesFind = $d.Content
$sevenDevicesFind.Find.Execute("7 devices, 6 weather types", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sevenDevicesPara = $sevenDevicesFind.Paragraphs(1)
$sevenDevicesStart = $sevenDevicesPara.Range.Start
$sevenDevicesEnd = $sevenDevicesPara.Range.End

# ---------------------------------------------------------------------------
# Step 2: Re-create the "_Hlk519930190" bookmark. It starts right before the
# "Device: All home appliances" run and runs down to the paragraph that now
# begins with "7 devices, 6 weather types...".
# ---------------------------------------------------------------------------
$deviceRange = $d.Content
$deviceRange.Find.Execute("Device: All home appliances", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$deviceRange.Collapse(1)
$hlk190Range = $d.Range($deviceRange.Start, $sevenDevicesStart)
$d.Bookmarks.Add("_Hlk519930190", $hlk190Range)

# ---------------------------------------------------------------------------
# Step 3: Add the "_Hlk519930332" bookmark. It starts right at the beginning
# of the "7 devices..." paragraph and covers the whole paragraph (i.e. runs
# through to the end of that same list item).
# ---------------------------------------------------------------------------
$hlk332Range = $d.Range($sevenDevicesStart, $sevenDevicesEnd)
$d.Bookmarks.Add("_Hlk519930332", $hlk332Range)

# ---------------------------------------------------------------------------
# Step 4: Split "Maximum Power Consumption Database" into "Maxi" / "mum Power
# Consumption Database" by dropping an (empty) "_GoBack" bookmark in between
# the two halves - this naturally splits the underlying run in two while
# keeping the original bold formatting on both pieces.
# ---------------------------------------------------------------------------
$titleRange = $d.Content
$titleRange.Find.Execute("Maximum Power Consumption Database", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $titleRange.Start + 4
$goBackRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $goBackRange)
